# Updated symbol list on Sat Feb 11 18:07:25 UTC 2023 with GitHub Actions
# Refreshes the crypto price/volume/hour snapshot cells in Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (kept as literal text, matching the
# original inlineStr cells, not auto-converted to number/percent).
$updates = [ordered]@{
    'D2' = '308.22'
    'E2' = '0.37%'
    'G2' = '18'
    'D3' = '40.91'
    'E3' = '-0.62%'
    'G3' = '18'
    'D4' = '5.122'
    'E4' = '1.45%'
    'G4' = '18'
    'D5' = '0.07625'
    'E5' = '0.14%'
    'G5' = '18'
    'D6' = '1.618'
    'E6' = '1.02%'
    'G6' = '18'
    'D7' = '2.474'
    'E7' = '2.27%'
    'G7' = '18'
    'D8' = '0.9092'
    'E8' = '0.47%'
    'G8' = '18'
    'D9' = '0.1259'
    'E9' = '29.48%'
    'G9' = '18'
    'D10' = '0.1817'
    'E10' = '2.84%'
    'G10' = '18'
    'D11' = '0.09122'
    'E11' = '-0.50%'
    'G11' = '18'
    'D12' = '0.04270'
    'E12' = '-1.47%'
    'G12' = '18'
    'D13' = '0.1044'
    'E13' = '-0.75%'
    'G13' = '18'
    'D14' = '0.001255'
    'E14' = '-0.47%'
    'G14' = '18'
    'D15' = '0.005775'
    'E15' = '-1.35%'
    'G15' = '18'
    'E16' = '-0.57%'
    'G16' = '18'
    'D17' = '4.284'
    'E17' = '0.79%'
    'G17' = '18'
    'E18' = '-0.65%'
    'G18' = '18'
    'D19' = '6.928'
    'E19' = '1.77%'
    'G19' = '18'
    'D20' = '0.1393'
    'E20' = '3.31%'
    'G20' = '18'
    'D21' = '0.2705'
    'E21' = '-0.68%'
    'G21' = '18'
    'D22' = '0.04041'
    'E22' = '-2.91%'
    'G22' = '18'
    'D23' = '0.001271'
    'E23' = '4.50%'
    'G23' = '18'
    'D24' = '0.004069'
    'E24' = '-0.04%'
    'G24' = '18'
    'E25' = '-2.23%'
    'G25' = '18'
    'E26' = '24.71%'
    'G26' = '18'
    'G27' = '18'
    'G28' = '18'
    'G29' = '18'
    'G30' = '18'
    'G31' = '18'
    'G32' = '18'
    'G33' = '18'
    'G34' = '18'
    'G35' = '18'
    'G36' = '18'
    'G37' = '18'
    'D38' = '0.02424'
    'E38' = '-0.42%'
    'G38' = '18'
    'D39' = '0.05231'
    'E39' = '1.59%'
    'G39' = '18'
    'D40' = '0.007833'
    'E40' = '-0.07%'
    'G40' = '18'
    'E41' = '-0.24%'
    'G41' = '18'
    'D42' = '0.006803'
    'E42' = '-3.78%'
    'G42' = '18'
    'E43' = '-0.85%'
    'G43' = '18'
    'D44' = '0.008060'
    'E44' = '-3.84%'
    'G44' = '18'
    'D45' = '0.3062'
    'E45' = '-8.29%'
    'G45' = '18'
    'D46' = '0.00006896'
    'E46' = '7.57%'
    'G46' = '18'
    'G47' = '18'
    'D48' = '0.1083'
    'E48' = '1,926.18%'
    'G48' = '18'
    'G49' = '18'
    'D50' = '0.00002103'
    'G50' = '18'
    'G51' = '18'
}

foreach ($ref in $updates.Keys) {
    $c = $ws.Range($ref)
    # Force text entry so values like "308.22" / "0.37%" / "18" stay as
    # strings instead of being auto-parsed into number/percent cells.
    $c.NumberFormat = "@"
    $c.Value = $updates[$ref]
    # Drop back to the default style so no stray number-format style
    # sticks to the cell (matches the original unstyled cells).
    $c.Style = "Normal"
}
